$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ablation text label change: "no" -> "el" (central-range visualization added)
$ws.Range("B2").Value = "el"

# Updated per-POI attribution values for the new ablation setting
$values = @(
    @{Row=3; Value=-0.01418817136436701},
    @{Row=4; Value=-0.01503446325659752},
    @{Row=5; Value=0.03324057534337044},
    @{Row=6; Value=0.02983036078512669},
    @{Row=7; Value=0.0356731005012989},
    @{Row=8; Value=-0.038100965321064},
    @{Row=9; Value=0.03669719398021698},
    @{Row=10; Value=0.03898508101701736},
    @{Row=11; Value=-0.06035224720835686},
    @{Row=12; Value=-0.02828722260892391},
    @{Row=13; Value=-0.01741178147494793},
    @{Row=14; Value=0.02714044786989689},
    @{Row=15; Value=0.008196081034839153},
    @{Row=16; Value=0.02671819180250168},
    @{Row=17; Value=0.004293262492865324},
    @{Row=18; Value=-0.0127369137480855},
    @{Row=19; Value=-0.08475329726934433},
    @{Row=20; Value=-0.05107626691460609},
    @{Row=21; Value=-0.01093430444598198},
    @{Row=22; Value=0.01788171194493771},
    @{Row=23; Value=0.0005421822424978018},
    @{Row=24; Value=-0.06557711213827133},
    @{Row=25; Value=0.01326756551861763},
    @{Row=26; Value=0.01878784783184528},
    @{Row=27; Value=-0.02680497989058495},
    @{Row=28; Value=-0.03166486695408821},
    @{Row=29; Value=0.05254664272069931},
    @{Row=30; Value=-0.004024737048894167},
    @{Row=31; Value=-0.03146779537200928},
    @{Row=32; Value=0.04504603520035744},
    @{Row=33; Value=0.05009873956441879},
    @{Row=34; Value=-0.06726761162281036},
    @{Row=35; Value=-0.0009220695937983692},
    @{Row=36; Value=0.004506095312535763},
    @{Row=37; Value=-0.02105730772018433},
    @{Row=38; Value=0.001260937424376607},
    @{Row=39; Value=-0.03853989392518997},
    @{Row=40; Value=0.0006460323347710073},
    @{Row=41; Value=-0.01117102708667517},
    @{Row=42; Value=0.05321109294891357},
    @{Row=43; Value=0.02866742201149464},
    @{Row=44; Value=-0.1117166131734848},
    @{Row=45; Value=-0.04144596308469772},
    @{Row=46; Value=0.03518811240792274},
    @{Row=47; Value=0.01006702240556479},
    @{Row=48; Value=-0.006136109586805105},
    @{Row=49; Value=-0.02939275465905666},
    @{Row=50; Value=-0.01789583824574947},
    @{Row=51; Value=-0.005569220054894686},
    @{Row=52; Value=-0.06661555916070938},
    @{Row=53; Value=-0.02575983107089996},
    @{Row=54; Value=-0.04909021779894829},
    @{Row=55; Value=0.01612738706171513},
    @{Row=56; Value=-0.05210376530885696},
    @{Row=57; Value=-0.05454245582222939},
    @{Row=58; Value=-0.001893242471851408},
    @{Row=59; Value=-0.04166665300726891},
    @{Row=60; Value=-0.03397499397397041},
    @{Row=61; Value=-0.08233810216188431},
    @{Row=62; Value=-0.06814239174127579},
    @{Row=63; Value=-0.01097240392118692},
    @{Row=64; Value=-0.05663706734776497},
    @{Row=65; Value=-0.02505508065223694},
    @{Row=66; Value=-0.0342939980328083},
    @{Row=67; Value=0.01235485635697842},
    @{Row=68; Value=-0.08740133792161942},
    @{Row=69; Value=0.005726086441427469},
    @{Row=70; Value=0.004742693156003952},
    @{Row=71; Value=-0.05000568926334381},
    @{Row=72; Value=0.01605942286550999},
    @{Row=73; Value=0.04612843319773674},
    @{Row=74; Value=-0.0157275628298521},
    @{Row=75; Value=-0.06928671151399612},
    @{Row=76; Value=0.02476442791521549},
    @{Row=77; Value=-0.02782525680959225},
    @{Row=78; Value=0.02152431383728981},
    @{Row=79; Value=-0.006752686575055122},
    @{Row=80; Value=-0.0006433501257561147},
    @{Row=81; Value=-0.04508950933814049},
    @{Row=82; Value=-0.02587642893195152},
    @{Row=83; Value=-0.0270231980830431},
    @{Row=84; Value=-0.006532173603773117},
    @{Row=85; Value=-0.05328363180160522},
    @{Row=86; Value=-0.09795290976762772},
    @{Row=87; Value=0.01896888576447964},
    @{Row=88; Value=-0.04846331849694252},
    @{Row=89; Value=-0.003385052317753434},
    @{Row=90; Value=-0.09989242255687714},
    @{Row=91; Value=-0.04493662342429161},
    @{Row=92; Value=0.007204334251582623},
    @{Row=93; Value=-0.01673361286520958},
    @{Row=94; Value=-0.06435196846723557},
    @{Row=95; Value=-0.07080777734518051},
    @{Row=96; Value=-0.1142769902944565},
    @{Row=97; Value=0.01272721122950315},
    @{Row=98; Value=-0.01295905280858278},
    @{Row=99; Value=0.003281347453594208},
    @{Row=100; Value=-0.05528610199689865},
    @{Row=101; Value=-0.07417339831590652},
    @{Row=102; Value=-0.007927651517093182},
    @{Row=103; Value=-0.08393402397632599},
    @{Row=104; Value=-0.0400700680911541},
    @{Row=105; Value=-0.02800003625452518},
    @{Row=106; Value=-0.0292271301150322},
    @{Row=107; Value=-0.06680292636156082},
    @{Row=108; Value=-0.008409772999584675},
    @{Row=109; Value=-0.02932162210345268},
    @{Row=110; Value=-0.08703675866127014},
    @{Row=111; Value=-0.09381962567567825},
    @{Row=112; Value=-0.06291677802801132},
    @{Row=113; Value=-0.005423332098871469},
    @{Row=114; Value=-0.01121558155864477},
    @{Row=115; Value=-0.1518655121326447},
    @{Row=116; Value=0.01616137847304344},
    @{Row=117; Value=0.01060223300009966},
    @{Row=118; Value=-0.02371115796267986},
    @{Row=119; Value=-0.03848818689584732},
    @{Row=120; Value=0.04758840426802635},
    @{Row=121; Value=-0.05567159503698349},
    @{Row=122; Value=-0.07417029142379761},
    @{Row=123; Value=0.06408835202455521},
    @{Row=124; Value=-0.05075591802597046},
    @{Row=125; Value=-0.03937947750091553},
    @{Row=126; Value=-0.105917327105999},
    @{Row=127; Value=-0.03317529708147049},
    @{Row=128; Value=-0.03965567052364349},
    @{Row=129; Value=-0.09042225033044815},
    @{Row=130; Value=-0.04746860265731812},
    @{Row=131; Value=-0.08238767087459564},
    @{Row=132; Value=-0.04472554102540016},
    @{Row=133; Value=0.006538607645779848},
    @{Row=134; Value=-0.0894000455737114},
    @{Row=135; Value=-0.08560368418693542},
    @{Row=136; Value=0.04131008312106133}
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 2).Value = $item.Value
}

Write-Host "Updated B2 and $($values.Count) value rows (B3:B136)."
